$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 72: "Vowels Game in a String"
$ws.Range("A72").Value = 3227
$ws.Range("B72").Value = "Vowels Game in a String"

# New row 73: "Maximum Number of Words You Can Type"
$ws.Range("A73").Value = 1935
$ws.Range("B73").Value = "Maximum Number of Words You Can Type"
$ws.Range("C73").Value = "Loop/String manipulation"
$ws.Range("D73").Value = "Loop in Loop"

# Fix typo in existing row 69 (E69): "string.fin(c)" -> "string.find(c)"
$ws.Range("E69").Value = "C++: vector, std::sort, string.find(c) != std::string::npos"

$ws.Range("E73").Value = "string stream >>"

# Update the active selection to match the author's final cursor position
$ws.Range("E73").Select()
